$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.02.2022 18:15"

# Row 7 (MOL Olomoucká): convert the delta price and date columns from
# text placeholders to real numeric values (AWS bash cmd line fix)
$ws.Range("D7").Value = 0.39
$ws.Range("E7").Value = 44613.75030092592
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
